# Update countries & provincias Spain
# Refresh the COVID-19 stats snapshot: update the "last updated" timestamp,
# refresh several countries' case numbers, and re-sort a handful of rows
# whose total-case counts changed enough to change their ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Junio de 2020 a las 05:12"

# Estados Unidos (row 4) - updated totals / recuperados
$ws.Cells.Item(4, 2).Value = 2462554
$ws.Cells.Item(4, 5).Value = 1297668

# Rows 60/61: Honduras overtakes Azerbaiyan in total cases, so they swap order
$ws.Cells.Item(60, 1).Value = "Honduras"
$ws.Cells.Item(60, 2).Value = 14571
$ws.Cells.Item(60, 3).Value = 628
$ws.Cells.Item(60, 4).Value = 1546
$ws.Cells.Item(60, 5).Value = 12608
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 12
$ws.Cells.Item(60, 8).Value = 417

$ws.Cells.Item(61, 1).Value = "Azerbaiyan"
$ws.Cells.Item(61, 2).Value = 14305
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 7768
$ws.Cells.Item(61, 5).Value = 6363
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 174

# Australia (row 74) - minor update to recuperados/muertes hoy/muertes
$ws.Cells.Item(74, 5).Value = 523
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 104

# Haiti (row 81) - updated totals
$ws.Cells.Item(81, 2).Value = 5429
$ws.Cells.Item(81, 3).Value = 105
$ws.Cells.Item(81, 4).Value = 512
$ws.Cells.Item(81, 5).Value = 4825
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 92

# Rows 147/148: Jamaica overtakes Libia, so they swap order
$ws.Cells.Item(147, 1).Value = "Jamaica"
$ws.Cells.Item(147, 2).Value = 678
$ws.Cells.Item(147, 3).Value = 8
$ws.Cells.Item(147, 4).Value = 521
$ws.Cells.Item(147, 5).Value = 147
$ws.Cells.Item(147, 8).Value = 10

$ws.Cells.Item(148, 1).Value = "Libia"
$ws.Cells.Item(148, 2).Value = 670
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 138
$ws.Cells.Item(148, 5).Value = 514
$ws.Cells.Item(148, 8).Value = 18

# Rows 200/201: Laos overtakes Santa Lucia, so they swap order
$ws.Cells.Item(200, 1).Value = "Laos"
$ws.Cells.Item(200, 4).Value = 19
$ws.Cells.Item(200, 5).Value = 0

$ws.Cells.Item(201, 1).Value = "Santa Lucia"
$ws.Cells.Item(201, 4).Value = 19
$ws.Cells.Item(201, 5).Value = 0

# Rows 206/207: Islas Turcas y Caicos overtakes San Cristobal y Nieves, so they swap order
$ws.Cells.Item(206, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(206, 2).Value = 15
$ws.Cells.Item(206, 3).Value = 1
$ws.Cells.Item(206, 4).Value = 11
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 8).Value = 1

$ws.Cells.Item(207, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(207, 2).Value = 15
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 15
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Rows 208/209: Islas Malvinas and Groenlandia swap display order (same data)
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 1).Value = "Groenlandia"
